$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.273.33"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.33"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.46"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6004"
$ws.Range("E6").Value = "  -4.19%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06981"
$ws.Range("E8").Value = "  -5.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.28"
$ws.Range("E10").Value = "  -6.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07616"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.08"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.754"
$ws.Range("E13").Value = "  -4.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6271"
$ws.Range("E14").Value = "  -6.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009670"
$ws.Range("E15").Value = "  -6.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.47"
$ws.Range("E16").Value = "  -3.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.858.10"
$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.648"
$ws.Range("E18").Value = "  -9.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.71"
$ws.Range("E19").Value = "  -5.74%  "

$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  -5.89%  "

$ws.Range("E22").Value = "  -5.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.37"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.964"
$ws.Range("E25").Value = "  -5.90%  "

$ws.Range("E26").Value = "  -4.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.54"
$ws.Range("E27").Value = "  -4.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.453"
$ws.Range("E28").Value = "  -2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06398"
$ws.Range("E29").Value = "  -12.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.440"
$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.838"
$ws.Range("E31").Value = "  -4.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.754"
$ws.Range("E32").Value = "  -6.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").Value = "  -5.76%  "

$ws.Range("E34").Value = "  -5.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6457"
$ws.Range("E35").Value = "  -9.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.542"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.745"
$ws.Range("E37").Value = "  -1.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01765"
$ws.Range("E38").Value = "  -3.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.589"
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.161.86"
$ws.Range("E40").Value = "  -5.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8928"
$ws.Range("E41").Value = "  -6.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.986.03"
$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.42"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.96"
$ws.Range("E45").Value = "  -4.93%  "

$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.587"
$ws.Range("E47").Value = "  -6.30%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.459"
$ws.Range("E48").Value = "  -5.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05548"
$ws.Range("E49").Value = "  -1.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4558"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.411"
$ws.Range("E51").Value = "  -7.72%  "
